{"js": "const replacements = [\n  [\"788\u00d73=\", \"835\u00d72=\"],\n  [\"357\u00d76=\", \"820\u00d73=\"],\n  [\"343\u00d72=\", \"986\u00d75=\"],\n  [\"935\u00d76=\", \"133\u00d77=\"],\n  [\"365\u00d72=\", \"706\u00d76=\"],\n  [\"393\u00d74=\", \"928\u00d75=\"],\n  [\"784\u00d76=\", \"375\u00d72=\"],\n  [\"296\u00d72=\", \"275\u00d74=\"],\n  [\"692\u00d78=\", \"522\u00d77=\"],\n  [\"371\u00d73=\", \"196\u00d76=\"],\n  [\"260\u00d75=\", \"687\u00d72=\"],\n  [\"979\u00d75=\", \"588\u00d76=\"],\n  [\"954\u00d78=\", \"530\u00d73=\"],\n  [\"824\u00d75=\", \"724\u00d75=\"],\n  [\"863\u00d77=\", \"740\u00d73=\"],\n  [\"718\u00d72=\", \"444\u00d72=\"],\n  [\"256\u00d75=\", \"594\u00d73=\"],\n  [\"831\u00d76=\", \"233\u00d76=\"],\n  [\"299\u00d79=\", \"367\u00d78=\"],\n  [\"371\u00d74=\", \"462\u00d77=\"],\n  [\"171\u00d78=\", \"523\u00d73=\"],\n  [\"122\u00d73=\", \"790\u00d73=\"],\n  [\"190\u00d76=\", \"104\u00d75=\"],\n  [\"643\u00d78=\", \"233\u00d74=\"],\n  [\"336\u00d73=\", \"104\u00d76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"788\u00d73=\", \"835\u00d72=\")\n    ,@(\"357\u00d76=\", \"820\u00d73=\")\n    ,@(\"343\u00d72=\", \"986\u00d75=\")\n    ,@(\"935\u00d76=\", \"133\u00d77=\")\n    ,@(\"365\u00d72=\", \"706\u00d76=\")\n    ,@(\"393\u00d74=\", \"928\u00d75=\")\n    ,@(\"784\u00d76=\", \"375\u00d72=\")\n    ,@(\"296\u00d72=\", \"275\u00d74=\")\n    ,@(\"692\u00d78=\", \"522\u00d77=\")\n    ,@(\"371\u00d73=\", \"196\u00d76=\")\n    ,@(\"260\u00d75=\", \"687\u00d72=\")\n    ,@(\"979\u00d75=\", \"588\u00d76=\")\n    ,@(\"954\u00d78=\", \"530\u00d73=\")\n    ,@(\"824\u00d75=\", \"724\u00d75=\")\n    ,@(\"863\u00d77=\", \"740\u00d73=\")\n    ,@(\"718\u00d72=\", \"444\u00d72=\")\n    ,@(\"256\u00d75=\", \"594\u00d73=\")\n    ,@(\"831\u00d76=\", \"233\u00d76=\")\n    ,@(\"299\u00d79=\", \"367\u00d78=\")\n    ,@(\"371\u00d74=\", \"462\u00d77=\")\n    ,@(\"171\u00d78=\", \"523\u00d73=\")\n    ,@(\"122\u00d73=\", \"790\u00d73=\")\n    ,@(\"190\u00d76=\", \"104\u00d75=\")\n    ,@(\"643\u00d78=\", \"233\u00d74=\")\n    ,@(\"336\u00d73=\", \"104\u00d76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
